$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'327.52"
$ws.Range("E2").Value = "'-0.82%"
$ws.Range("D3").Value = "'44.06"
$ws.Range("E3").Value = "'1.22%"
$ws.Range("D4").Value = "'5.554"
$ws.Range("E4").Value = "'-0.78%"
$ws.Range("D5").Value = "'0.08047"
$ws.Range("E5").Value = "'-1.78%"
$ws.Range("D6").Value = "'1.903"
$ws.Range("E6").Value = "'-0.11%"
$ws.Range("D7").Value = "'4.273"
$ws.Range("E7").Value = "'-3.09%"
$ws.Range("D8").Value = "'2.559"
$ws.Range("E8").Value = "'-9.11%"
$ws.Range("D9").Value = "'0.9441"
$ws.Range("E9").Value = "'0.13%"
$ws.Range("D10").Value = "'0.1165"
$ws.Range("E10").Value = "'-2.75%"
$ws.Range("E11").Value = "'-4.36%"
$ws.Range("D12").Value = "'0.09671"
$ws.Range("E12").Value = "'-2.62%"
$ws.Range("D13").Value = "'0.04381"
$ws.Range("E13").Value = "'0.57%"
$ws.Range("D14").Value = "'0.1068"
$ws.Range("E14").Value = "'-0.11%"
$ws.Range("D15").Value = "'0.001275"
$ws.Range("E15").Value = "'-0.08%"
$ws.Range("D16").Value = "'0.005994"
$ws.Range("E16").Value = "'0.12%"
$ws.Range("B17").Value = "HotbitToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D17").Value = "'0.004279"
$ws.Range("E17").Value = "'-0.52%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.405"
$ws.Range("E18").Value = "'-2.77%"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "'0.3491"
$ws.Range("E19").Value = "'-1.30%"
$ws.Range("B20").Value = "MCDex"
$ws.Range("C20").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D20").Value = "'9.941"
$ws.Range("E20").Value = "'13.77%"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "'0.1379"
$ws.Range("E21").Value = "'0.69%"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").Value = "'0.2508"
$ws.Range("E22").Value = "'-0.55%"
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").Value = "'0.04198"
$ws.Range("E23").Value = "'-4.64%"
$ws.Range("B24").Value = "BitKan"
$ws.Range("C24").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D24").Value = "'0.001247"
$ws.Range("E24").Value = "'0.54%"
$ws.Range("E25").Value = "'2.29%"
$ws.Range("D26").Value = "'0.0003996"
$ws.Range("E26").Value = "'-0.25%"
$ws.Range("D38").Value = "'0.02643"
$ws.Range("E38").Value = "'-6.40%"
$ws.Range("D39").Value = "'0.05478"
$ws.Range("E39").Value = "'-4.54%"
$ws.Range("D40").Value = "'0.007567"
$ws.Range("E40").Value = "'-4.16%"
$ws.Range("E41").Value = "'-1.81%"
$ws.Range("D42").Value = "'0.007996"
$ws.Range("E42").Value = "'-18.46%"
$ws.Range("D43").Value = "'0.002011"
$ws.Range("E43").Value = "'-4.53%"
$ws.Range("D44").Value = "'0.008841"
$ws.Range("E44").Value = "'-11.39%"
$ws.Range("D45").Value = "'0.00006936"
$ws.Range("E45").Value = "'-5.21%"
$ws.Range("E46").Value = "'-0.10%"
$ws.Range("D47").Value = "'0.002274"
$ws.Range("E47").Value = "'-0.25%"
$ws.Range("D48").Value = "'0.005881"
$ws.Range("E48").Value = "'67.97%"
$ws.Range("D49").Value = "'0.00002104"
$ws.Range("E49").Value = "'-0.10%"
$ws.Range("D50").Value = "'0.0002003"
$ws.Range("E50").Value = "'-0.10%"

# Reset style (quotePrefix) introduced by the leading apostrophe trick back to Normal
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Style = "Normal"
